# Re-upload of the manual testing workbook following amendments.
# - search_function row (C13): append note about swapping the gspread-specific
#   exception for a general Exception, keeping the existing "TRY/EXCEPT:" bold run.
# - main_program_call row (C15): append note about the bug found before submission
#   (the search-branch loop / broad except fix) and how it was resolved.
# - Row heights grow to fit the extra wrapped lines of text.
# - Selection/scroll moves down to the newly-edited row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C13: search_function comments -> append the exception-handling note ---
$c13 = $ws.Range("C13")
$c13Old = $c13.Value()
$c13Addition = "Had to replace the gspread exception I wrote in with a general exception as anything more specific throws a bug, this should not affect the outcome of the code."
$c13New = $c13Old + "`n" + $c13Addition
$c13.Value = $c13New

# Re-apply the bold formatting to the "TRY/EXCEPT:" label that gets reset when
# the whole cell value is overwritten.
$boldIdx = $c13New.IndexOf("TRY/EXCEPT:")
if ($boldIdx -ge 0) {
    $c13.Characters($boldIdx + 1, [string]"TRY/EXCEPT:".Length).Font.Bold = $true
}

# --- C15: main_program_call comments -> append the late bug-fix note ---
$c15 = $ws.Range("C15")
$c15Old = $c15.Value()
$c15Addition = "Ran into a fairly significant bug before submission due to the except statement in the search function not triggering. I changed the exception to a general Exception as any other was continually throwing a bug. The loop in the search branch of this function was flawed and was not allowing clear looping of the search functionality from unsuccessful and successful search results. The code is now functional."
$c15New = $c15Old + "`n" + $c15Addition
$c15.Value = $c15New

# --- Row heights grow to accommodate the extra wrapped lines ---
$ws.Rows(13).RowHeight = 375
$ws.Rows(15).RowHeight = 409.5

# --- Move the active selection down to the row that was just edited ---
$ws.Range("C15").Select()
